$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 31; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()

    if ($a -eq 16) { $ws.Cells.Item($r, 1).Value = 12 }
    elseif ($a -eq -16) { $ws.Cells.Item($r, 1).Value = -12 }
    elseif ($a -eq 11.314) { $ws.Cells.Item($r, 1).Value = 8.485 }
    elseif ($a -eq -11.314) { $ws.Cells.Item($r, 1).Value = -8.485 }

    if ($b -eq 16) { $ws.Cells.Item($r, 2).Value = 12 }
    elseif ($b -eq -16) { $ws.Cells.Item($r, 2).Value = -12 }
    elseif ($b -eq 11.314) { $ws.Cells.Item($r, 2).Value = 8.485 }
    elseif ($b -eq -11.314) { $ws.Cells.Item($r, 2).Value = -8.485 }
}

$ws.Range("F22").Select()
